$d = $word.ActiveDocument

# 1. "Algorithm One (G):" -> "Algorithm One (G, exclude):"
$d.Content.Find.Execute(" (G):", $false, $false, $false, $false, $false, $true, 1, $false, " (G, exclude):", 2)

# 2. "Algorithm Two(G):" -> "Algorithm Two(G, exclude):"
$d.Content.Find.Execute("(G):", $false, $false, $false, $false, $false, $true, 1, $false, "(G, exclude):", 2)

# 3. Both "Input: an undirected, weighted, connected graph G" lines get the exclude note appended
$d.Content.Find.Execute("Input: an undirected, weighted, connected graph G", $false, $false, $false, $false, $false, $true, 1, $false, "Input: an undirected, weighted, connected graph G and a Node exclude, if needed, otherwise it is null", 2)

# 4. Both "< minimum_distance " conditions get "and v is not exclude" appended
$d.Content.Find.Execute(" < minimum_distance ", $false, $false, $false, $false, $false, $true, 1, $false, " < minimum_distance and v is not exclude ", 2)
